$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: J0 "Objectifs" paragraph - append new sentence about deciding to
# write down hours to help write the real planning later.
# ---------------------------------------------------------------------------
$objParaRange = $d.Paragraphs(4).Range
[void]$objParaRange.InsertAfter(" J’ai ")
$objParaRange = $d.Paragraphs(4).Range
[void]$objParaRange.InsertAfter("décidé")
$objParaRange = $d.Paragraphs(4).Range
[void]$objParaRange.InsertAfter(" d’écrire les heures pour m’aider ")
$objParaRange = $d.Paragraphs(4).Range
[void]$objParaRange.InsertAfter("à")
$objParaRange = $d.Paragraphs(4).Range
[void]$objParaRange.InsertAfter(" rédiger le planning réel par la suite")

# ---------------------------------------------------------------------------
# Change 2: J1 "Objectifs" paragraph - merge the first two runs ("Aujourd'hui
# est le deuxième jour du TPI" + ", je vais m'occuper ") into a single run,
# leaving the following runs ("de la gestion" / " des utilisateurs et des
# groupes" / ".") untouched.
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("Aujourd") -and $t.Contains("deuxième")) {
        $j1ObjIndex = $i
        $found = $true
        break
    }
}
$j1ObjRange = $d.Paragraphs($j1ObjIndex).Range
$xml = "<w:p $wns><w:r><w:t xml:space=`"preserve`">Aujourd’hui est le deuxième jour du TPI, je vais m’occuper </w:t></w:r><w:r><w:t>de la gestion</w:t></w:r><w:r><w:t xml:space=`"preserve`"> des utilisateurs et des groupes</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"
[void]$j1ObjRange.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change 3: rework the 08:30 / 09:00 block and add the rest of the day's log.
# ---------------------------------------------------------------------------
# Locate the "08:30" paragraph freshly (index did not move from the edits
# above, since none of them added/removed paragraphs).
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("08:30")) {
        $idx0830 = $i
        $found = $true
        break
    }
}

# Insert a brand-new paragraph before "08:30" asking how to store the salt.
$p0830 = $d.Paragraphs($idx0830).Range
[void]$p0830.InsertParagraphBefore()
$idxSalt = $idx0830
$idx0830 = $idx0830 + 1

$saltXml = "<w:p $wns><w:r><w:t xml:space=`"preserve`">Je me demande de quelle manière je vais stocker le </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>salt</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> dans la base de donnée</w:t></w:r></w:p>"
[void]$d.Paragraphs($idxSalt).Range.InsertXML($saltXml)

# Rewrite the "08:30" paragraph with the newly found info, followed by the
# original sentence about finishing the registration.
$xml0830 = "<w:p $wns><w:r><w:t>08:30 :</w:t></w:r><w:r><w:t xml:space=`"preserve`"> J’ai trouvé une manière sécurisée pour stocker le </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>salt</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> dans la base de donnée</w:t></w:r><w:r><w:t xml:space=`"preserve`"> Je fini l’enregistrement. Je commence la connexion.</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx0830).Range.InsertXML($xml0830)

# The "09:00" paragraph follows right after; locate it and strip the
# _GoBack bookmark that currently lives inside it (it will move to the new
# final paragraph of the day).
$idx0900 = $idx0830 + 1
$xml0900 = "<w:p $wns><w:r><w:t>09:00 : Je fini la connexion. Je commence la modification de mot de passe de l’utilisateur.</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx0900).Range.InsertXML($xml0900)

# Append the rest of the day's log as new paragraphs after "09:00".
$p0900 = $d.Paragraphs($idx0900).Range
[void]$p0900.InsertParagraphAfter()
$idx0920 = $idx0900 + 1
$xml0920 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">09:20 : Mon camarade </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Yvelin</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> RAY me demande les normes de nommages de base de donnée de l’école</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx0920).Range.InsertXML($xml0920)

$p0920 = $d.Paragraphs($idx0920).Range
[void]$p0920.InsertParagraphAfter()
$idx0925 = $idx0920 + 1
$xml0925 = "<w:p $wns><w:r><w:t>09:25 : Je reprends mon travail.</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx0925).Range.InsertXML($xml0925)

$p0925 = $d.Paragraphs($idx0925).Range
[void]$p0925.InsertParagraphAfter()
$idx0935 = $idx0925 + 1
$xml0935 = "<w:p $wns><w:r><w:t>09:35</w:t></w:r><w:r><w:t xml:space=`"preserve`"> : Anthony me montre comment mettre un </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>footer</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> en bas de page</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx0935).Range.InsertXML($xml0935)

$p0935 = $d.Paragraphs($idx0935).Range
[void]$p0935.InsertParagraphAfter()
$idx0940 = $idx0935 + 1
$xml0940 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">09:40 : Je </w:t></w:r><w:r><w:t>prends</w:t></w:r><w:r><w:t xml:space=`"preserve`"> ma pause</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx0940).Range.InsertXML($xml0940)

$p0940 = $d.Paragraphs($idx0940).Range
[void]$p0940.InsertParagraphAfter()
$idx1005 = $idx0940 + 1
$xml1005 = "<w:p $wns><w:r><w:t>10:05 : Je reprends la modification du profil utilisateur</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx1005).Range.InsertXML($xml1005)

$p1005 = $d.Paragraphs($idx1005).Range
[void]$p1005.InsertParagraphAfter()
$idx1030 = $idx1005 + 1
$xml1030 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">10:30 : </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Costantino</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> VOLTA m’a aidé à corriger des fautes de français sur le site</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx1030).Range.InsertXML($xml1030)

$p1030 = $d.Paragraphs($idx1030).Range
[void]$p1030.InsertParagraphAfter()
$idx1035 = $idx1030 + 1
$xml1035 = "<w:p $wns><w:r><w:t>10:35 : J’ai fini la gestion utilisateur sur le site. Je vais avancer la documentation sur la partie de gestion utilisateur.</w:t></w:r></w:p>"
[void]$d.Paragraphs($idx1035).Range.InsertXML($xml1035)

# Re-add the _GoBack bookmark at the end of the new, final paragraph.
$p1035 = $d.Paragraphs($idx1035).Range
$p1035.Collapse(0)
[void]$p1035.Bookmarks.Add("_GoBack")

# Finally, remove the trailing empty paragraph that used to close the
# section (there is one extra empty paragraph right before the sectPr now).
$lastIdx = $d.Paragraphs.Count
$lastText = $d.Paragraphs($lastIdx).Range.Text
if ($lastText -eq "") {
    [void]$d.Paragraphs($lastIdx).Range.Delete()
}

Write-Output "edit complete"
